$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 772.2941
$ws.Range("I19").Value = 633.7692
$ws.Range("J19").Value = 1222.5
$ws.Range("K19").Value = 633.7692
$ws.Range("L19").Value = 1222.5
$ws.Range("M19").Value = -458.7692
$ws.Range("N19").Value = -1572.5
$ws.Range("H70").Value = 2750
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 2750
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -10872
$ws.Range("H80").Value = 766.53845
$ws.Range("I80").Value = 390.5
$ws.Range("J80").Value = 933.6667
$ws.Range("K80").Value = 1171.5
$ws.Range("L80").Value = 2801.0001
$ws.Range("M80").Value = -173.5
$ws.Range("N80").Value = -4797.0001
$ws.Range("H83").Value = 766.53845
$ws.Range("I83").Value = 390.5
$ws.Range("J83").Value = 933.6667
$ws.Range("K83").Value = 3514.5
$ws.Range("L83").Value = 8403.0003
$ws.Range("M83").Value = 1477.5
$ws.Range("N83").Value = -18387.0003
$ws.Range("H129").Value = 881.04
$ws.Range("I129").Value = 513.5
$ws.Range("J129").Value = 972.925
$ws.Range("K129").Value = 1540.5
$ws.Range("L129").Value = 2918.775
$ws.Range("M129").Value = 3459.5
$ws.Range("N129").Value = -12918.775
$ws.Range("H138").Value = 1724.54
$ws.Range("I138").Value = 533.41174
$ws.Range("J138").Value = 1968.506
$ws.Range("K138").Value = 1600.23522
$ws.Range("L138").Value = 5905.518
$ws.Range("M138").Value = 3539.76478
$ws.Range("N138").Value = -16185.518
$ws.Range("H141").Value = 2007.1111
$ws.Range("I141").Value = 2007.1111
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6021.3333
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -841.3333000000002
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 77078104
$ws.Range("I61").Value = 100101200
$ws.Range("J61").Value = 334466.66
$ws.Range("K61").Value = 100101200
$ws.Range("L61").Value = 334466.66
$ws.Range("M61").Value = -100100988
$ws.Range("N61").Value = -334890.66
$ws.Range("H136").Value = 77078104
$ws.Range("I136").Value = 100101200
$ws.Range("J136").Value = 334466.66
$ws.Range("K136").Value = 300303600
$ws.Range("L136").Value = 1003399.98
$ws.Range("M136").Value = -300301050
$ws.Range("N136").Value = -1008499.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1023
$ws.Range("I20").Value = 930.2857
$ws.Range("J20").Value = 1239.3334
$ws.Range("K20").Value = 930.2857
$ws.Range("L20").Value = 1239.3334
$ws.Range("M20").Value = -683.2857
$ws.Range("N20").Value = -1733.3334
$ws.Range("H105").Value = 378112.5
$ws.Range("I105").Value = 503316.66
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 503316.66
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -501569.66
$ws.Range("N105").Value = -5994
$ws.Range("H107").Value = 2546.6365
$ws.Range("I107").Value = 2500.5
$ws.Range("J107").Value = 2556.889
$ws.Range("K107").Value = 2500.5
$ws.Range("L107").Value = 2556.889
$ws.Range("M107").Value = -580.5
$ws.Range("N107").Value = -6396.889
$ws.Range("H126").Value = 75000
$ws.Range("J126").Value = 75000
$ws.Range("L126").Value = 75000
$ws.Range("N126").Value = -84880
$ws.Range("H134").Value = 4678.4707
$ws.Range("I134").Value = 3829.9312
$ws.Range("J134").Value = 9600
$ws.Range("K134").Value = 11489.7936
$ws.Range("L134").Value = 28800
$ws.Range("M134").Value = -8954.793600000001
$ws.Range("N134").Value = -33870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 411.68292
$ws.Range("I22").Value = 435.12903
$ws.Range("J22").Value = 339
$ws.Range("K22").Value = 435.12903
$ws.Range("L22").Value = 339
$ws.Range("M22").Value = -85.12903
$ws.Range("N22").Value = -1039
$ws.Range("H31").Value = 2816.65
$ws.Range("I31").Value = 1374.5238
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 1374.5238
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -1079.5238
$ws.Range("N31").Value = -3790
$ws.Range("H34").Value = 2816.65
$ws.Range("I34").Value = 1374.5238
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 1374.5238
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -1172.5238
$ws.Range("N34").Value = -3604
$ws.Range("H58").Value = 34486360
$ws.Range("I58").Value = 38465324
$ws.Range("K58").Value = 38465324
$ws.Range("M58").Value = -38465121
$ws.Range("H99").Value = 2873.3333
$ws.Range("I99").Value = 3084.6155
$ws.Range("K99").Value = 3084.6155
$ws.Range("M99").Value = -1586.6155
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 1220.7333
$ws.Range("I107").Value = 490.1111
$ws.Range("J107").Value = 2316.6667
$ws.Range("K107").Value = 490.1111
$ws.Range("L107").Value = 2316.6667
$ws.Range("M107").Value = 1429.8889
$ws.Range("N107").Value = -6156.6667
$ws.Range("H126").Value = 2873.3333
$ws.Range("I126").Value = 3084.6155
$ws.Range("K126").Value = 9253.8465
$ws.Range("M126").Value = -6783.8465
$ws.Range("H136").Value = 34486360
$ws.Range("I136").Value = 38465324
$ws.Range("K136").Value = 115395972
$ws.Range("M136").Value = -115393422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5385.2856
$ws.Range("I70").Value = 5221.5386
$ws.Range("J70").Value = 5651.375
$ws.Range("K70").Value = 5221.5386
$ws.Range("L70").Value = 5651.375
$ws.Range("M70").Value = -4951.5386
$ws.Range("N70").Value = -6191.375
$ws.Range("H73").Value = 5385.2856
$ws.Range("I73").Value = 5221.5386
$ws.Range("J73").Value = 5651.375
$ws.Range("K73").Value = 5221.5386
$ws.Range("L73").Value = 5651.375
$ws.Range("M73").Value = -4285.5386
$ws.Range("N73").Value = -7523.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4554.364
$ws.Range("I40").Value = 4233.3335
$ws.Range("J40").Value = 5999
$ws.Range("K40").Value = 4233.3335
$ws.Range("L40").Value = 5999
$ws.Range("M40").Value = -4097.3335
$ws.Range("N40").Value = -6271
$ws.Range("H122").Value = 3227.5
$ws.Range("I122").Value = 2926.1538
$ws.Range("K122").Value = 8778.4614
$ws.Range("M122").Value = -6328.4614
$ws.Range("H136").Value = 211642.86
$ws.Range("I136").Value = 170833.33
$ws.Range("J136").Value = 242250
$ws.Range("K136").Value = 512499.99
$ws.Range("L136").Value = 726750
$ws.Range("M136").Value = -509949.99
$ws.Range("N136").Value = -731850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 63395.812
$ws.Range("I100").Value = 72185.71000000001
$ws.Range("K100").Value = 144371.42
$ws.Range("M100").Value = -143830.42
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 2550
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -15440
